$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '39.868.97'
$ws.Cells.Item(2, 5).Value = '  +1.13%  '
$ws.Cells.Item(3, 4).Value = '2.196.73'
$ws.Cells.Item(3, 5).Value = '  +1.72%  '
$ws.Cells.Item(4, 5).Value = '  -0.08%  '
$ws.Cells.Item(5, 4).Formula = '''227.84'
$ws.Cells.Item(5, 5).Value = '  -0.62%  '
$ws.Cells.Item(6, 5).Value = '  +0.93%  '
$ws.Cells.Item(7, 4).Formula = '''63.50'
$ws.Cells.Item(7, 5).Value = '  +0.88%  '
$ws.Cells.Item(8, 5).Value = '  +0.01%  '
$ws.Cells.Item(9, 4).Formula = '''0.394'
$ws.Cells.Item(9, 5).Value = '  -0.43%  '
$ws.Cells.Item(10, 5).Value = '  -0.76%  '
$ws.Cells.Item(11, 5).Value = '  +0.69%  '
$ws.Cells.Item(12, 4).Formula = '''16.06'
$ws.Cells.Item(12, 5).Value = '  +0.14%  '
$ws.Cells.Item(13, 4).Value = '2.521.21'
$ws.Cells.Item(13, 5).Value = '  +1.34%  '
$ws.Cells.Item(14, 4).Formula = '''22.04'
$ws.Cells.Item(14, 5).Value = '  -0.72%  '
$ws.Cells.Item(15, 4).Formula = '''0.819'
$ws.Cells.Item(15, 5).Value = '  +0.12%  '
$ws.Cells.Item(16, 4).Formula = '''5.56'
$ws.Cells.Item(16, 5).Value = '  -0.19%  '
$ws.Cells.Item(17, 4).Value = '2.197.18'
$ws.Cells.Item(17, 5).Value = '  +1.18%  '
$ws.Cells.Item(18, 4).Value = '39.888.09'
$ws.Cells.Item(18, 5).Value = '  +1.19%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0922'
$ws.Cells.Item(19, 5).Value = '  +8.23%  '
$ws.Cells.Item(20, 4).Formula = '''72.03'
$ws.Cells.Item(20, 5).Value = '  -0.43%  '
$ws.Cells.Item(21, 4).Formula = '''6.05'
$ws.Cells.Item(21, 5).Value = '  -1.50%  '
$ws.Cells.Item(22, 4).Formula = '''231.10'
$ws.Cells.Item(22, 5).Value = '  +1.17%  '
$ws.Cells.Item(23, 5).Value = '  +0.05%  '
$ws.Cells.Item(24, 5).Value = '  +0.09%  '
$ws.Cells.Item(25, 4).Formula = '''2.37'
$ws.Cells.Item(25, 5).Value = '  -0.26%  '
$ws.Cells.Item(26, 4).Formula = '''170.98'
$ws.Cells.Item(26, 5).Value = '  -0.69%  '
$ws.Cells.Item(27, 4).Formula = '''9.52'
$ws.Cells.Item(27, 5).Value = '  -2.55%  '
$ws.Cells.Item(28, 5).Value = '  +1.40%  '
$ws.Cells.Item(29, 4).Formula = '''1.47'
$ws.Cells.Item(29, 5).Value = '  +3.54%  '
$ws.Cells.Item(30, 4).Formula = '''19.99'
$ws.Cells.Item(30, 5).Value = '  +1.76%  '
$ws.Cells.Item(31, 4).Formula = '''2.69'
$ws.Cells.Item(31, 5).Value = '  +4.24%  '
$ws.Cells.Item(32, 4).Formula = '''0.123'
$ws.Cells.Item(32, 5).Value = '  +1.00%  '
$ws.Cells.Item(33, 4).Formula = '''4.55'
$ws.Cells.Item(33, 5).Value = '  -2.16%  '
$ws.Cells.Item(34, 5).Value = '  -2.03%  '
$ws.Cells.Item(35, 4).Formula = '''7.03'
$ws.Cells.Item(35, 5).Value = '  +0.40%  '
$ws.Cells.Item(36, 2).Value = 'Hedera'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(36, 4).Formula = '''0.0622'
$ws.Cells.Item(36, 5).Value = '  +0.08%  '
$ws.Cells.Item(37, 2).Value = 'RenderToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(37, 4).Formula = '''3.87'
$ws.Cells.Item(37, 5).Value = '  +9.31%  '
$ws.Cells.Item(38, 4).Formula = '''2.44'
$ws.Cells.Item(38, 5).Value = '  +0.47%  '
$ws.Cells.Item(39, 5).Value = '  -0.09%  '
$ws.Cells.Item(40, 4).Formula = '''4.99'
$ws.Cells.Item(40, 5).Value = '  +17.42%  '
$ws.Cells.Item(41, 4).Formula = '''103.23'
$ws.Cells.Item(41, 5).Value = '  -0.39%  '
$ws.Cells.Item(42, 4).Formula = '''0.0229'
$ws.Cells.Item(42, 5).Value = '  -0.81%  '
$ws.Cells.Item(43, 4).Formula = '''17.88'
$ws.Cells.Item(43, 5).Value = '  -1.06%  '
$ws.Cells.Item(44, 2).Value = 'Maker'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(44, 4).Value = '1.516.12'
$ws.Cells.Item(44, 5).Value = '  -0.78%  '
$ws.Cells.Item(45, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(45, 4).Formula = '''1.22'
$ws.Cells.Item(45, 5).Value = '  +3.18%  '
$ws.Cells.Item(46, 4).Formula = '''7.96'
$ws.Cells.Item(46, 5).Value = '  +2.36%  '
$ws.Cells.Item(47, 4).Formula = '''0.0924'
$ws.Cells.Item(47, 5).Value = '  -0.75%  '
$ws.Cells.Item(48, 5).Value = '  -0.22%  '
$ws.Cells.Item(49, 5).Value = '  -0.36%  '
$ws.Cells.Item(50, 4).Formula = '''0.000192'
$ws.Cells.Item(50, 5).Value = '  +32.05%  '
$ws.Cells.Item(51, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(51, 4).Value = '2.400.79'
$ws.Cells.Item(51, 5).Value = '  +1.24%  '
